# Apply scheduled-update refresh for Linea 141 (scrape timestamp 04:03:00)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" : full refresh of the arrivals table (rows 6-13)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:03:00"

$ws1.Range("A6").Value = "04:03:00"
$ws1.Range("B6").Value = "04:03"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 0

$ws1.Range("A7").Value = "04:03:00"
$ws1.Range("B7").Value = "04:46"
$ws1.Range("C7").Value = "215A_EL PATO"
$ws1.Range("D7").Value = 43

$ws1.Range("A8").Value = "04:03:00"
$ws1.Range("B8").Value = "04:53"
$ws1.Range("C8").Value = "11_ETCHEVERRY"
$ws1.Range("D8").Value = 50

$ws1.Range("A9").Value = "04:03:00"
$ws1.Range("B9").Value = "05:16"
$ws1.Range("C9").Value = "17_ROMERO"
$ws1.Range("D9").Value = 73

$ws1.Range("A10").Value = "04:03:00"
$ws1.Range("B10").Value = "05:22"
$ws1.Range("C10").Value = "23_HERNANDEZ"
$ws1.Range("D10").Value = 79

$ws1.Range("A11").Value = "04:03:00"
$ws1.Range("B11").Value = "05:35"
$ws1.Range("C11").Value = "215B_EL PATO"
$ws1.Range("D11").Value = 92

$ws1.Range("A12").Value = "04:03:00"
$ws1.Range("B12").Value = "05:41"
$ws1.Range("C12").Value = "14_ABASTO"
$ws1.Range("D12").Value = 98

$ws1.Range("A13").Value = "04:03:00"
$ws1.Range("B13").Value = "05:46"
$ws1.Range("C13").Value = "15_ABASTO"
$ws1.Range("D13").Value = 103

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" : refresh timestamps + updated wait minutes
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:03:00"

$ws2.Range("A6").Value = "04:03:00"
$ws2.Range("D6").Value = 43

$ws2.Range("A7").Value = "04:03:00"
$ws2.Range("B7").Value = "05:35"
$ws2.Range("D7").Value = 92

# ---------------------------------------------------------------------------
# Sheet "6203-6173" : first scraped row ever recorded for this stop pair
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:03:00"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"

$ws3.Range("A6").Value = "04:03:00"
$ws3.Range("B6").Value = "05:44"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 101
$ws3.Range("E6").Value = "L6173"
